$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '51.272.38'
$c.Style = "Normal"

$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +4.57%  '
$c.Style = "Normal"

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.743.29'
$c.Style = "Normal"

$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  +4.43%  '
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '115.39'
$c.Style = "Normal"

$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +3.65%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '331.96'
$c.Style = "Normal"

$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +2.93%  '
$c.Style = "Normal"

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.538'
$c.Style = "Normal"

$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +2.39%  '
$c.Style = "Normal"

$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  -0.01%  '
$c.Style = "Normal"

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.572'
$c.Style = "Normal"

$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  +5.70%  '
$c.Style = "Normal"

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '41.42'
$c.Style = "Normal"

$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +4.51%  '
$c.Style = "Normal"

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0854'
$c.Style = "Normal"

$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +5.52%  '
$c.Style = "Normal"

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '20.19'
$c.Style = "Normal"

$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  +2.40%  '
$c.Style = "Normal"

$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  +2.32%  '
$c.Style = "Normal"

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '7.60'
$c.Style = "Normal"

$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +5.03%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '3.182.24'
$c.Style = "Normal"

$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  +4.82%  '
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '2.777.62'
$c.Style = "Normal"

$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +5.43%  '
$c.Style = "Normal"

$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +2.91%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '51.361.43'
$c.Style = "Normal"

$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  +4.84%  '
$c.Style = "Normal"

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '3.21'
$c.Style = "Normal"

$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +6.44%  '
$c.Style = "Normal"

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '13.42'
$c.Style = "Normal"

$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  +3.96%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.83'
$c.Style = "Normal"

$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +2.16%  '
$c.Style = "Normal"

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.0₃0972'
$c.Style = "Normal"

$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  +3.17%  '
$c.Style = "Normal"

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '277.55'
$c.Style = "Normal"

$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +3.24%  '
$c.Style = "Normal"

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '69.45'
$c.Style = "Normal"

$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  +1.38%  '
$c.Style = "Normal"

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.66'
$c.Style = "Normal"

$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +5.04%  '
$c.Style = "Normal"

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '26.72'
$c.Style = "Normal"

$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +2.37%  '
$c.Style = "Normal"

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"

$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '10.15'
$c.Style = "Normal"

$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -0.22%  '
$c.Style = "Normal"

$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  -0.68%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.140'
$c.Style = "Normal"

$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  +1.42%  '
$c.Style = "Normal"

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '35.02'
$c.Style = "Normal"

$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c.Style = "Normal"

$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  +0.85%  '
$c.Style = "Normal"

$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +1.53%  '
$c.Style = "Normal"

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.0822'
$c.Style = "Normal"

$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +2.74%  '
$c.Style = "Normal"

$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c.Style = "Normal"

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '19.02'
$c.Style = "Normal"

$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c.Style = "Normal"

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.08'
$c.Style = "Normal"

$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  +2.35%  '
$c.Style = "Normal"

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '4.98'
$c.Style = "Normal"

$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  +0.72%  '
$c.Style = "Normal"

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '3.24'
$c.Style = "Normal"

$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  +4.08%  '
$c.Style = "Normal"

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '126.78'
$c.Style = "Normal"

$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -0.17%  '
$c.Style = "Normal"

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '23.28'
$c.Style = "Normal"

$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +5.18%  '
$c.Style = "Normal"

$c = $ws.Range('B42')
$c.NumberFormat = "@"
$c.Value = 'VeChain'
$c.Style = "Normal"

$c = $ws.Range('C42')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c.Style = "Normal"

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.0347'
$c.Style = "Normal"

$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  +9.51%  '
$c.Style = "Normal"

$c = $ws.Range('B43')
$c.NumberFormat = "@"
$c.Value = 'WEMIXToken'
$c.Style = "Normal"

$c = $ws.Range('C43')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c.Style = "Normal"

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.30'
$c.Style = "Normal"

$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  +7.91%  '
$c.Style = "Normal"

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.113'
$c.Style = "Normal"

$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  +2.72%  '
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.43'
$c.Style = "Normal"

$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  +13.02%  '
$c.Style = "Normal"

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.088.51'
$c.Style = "Normal"

$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  +1.27%  '
$c.Style = "Normal"

$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  +3.57%  '
$c.Style = "Normal"

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '5.50'
$c.Style = "Normal"

$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  +6.00%  '
$c.Style = "Normal"

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '8.96'
$c.Style = "Normal"

$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  +0.90%  '
$c.Style = "Normal"

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '59.77'
$c.Style = "Normal"

$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  +2.10%  '
$c.Style = "Normal"
